# Update the "cryptos" price/volume table with refreshed figures.
# (Mirrors a scheduled GitHub Actions data refresh.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several Price cells (column D) hold values that look like plain decimals
# (e.g. "213.10"). Column D has no explicit number format (General), so
# assigning such a string would be auto-parsed into a Double and would lose
# the trailing zero / exact text. Force those specific cells to Text first,
# write the literal string, then restore the default "Normal" style so the
# cell formatting itself is left unchanged (matches the original workbook).
$forceTextCells = @(
    'D5', 'D6', 'D9', 'D10', 'D11', 'D14', 'D15', 'D17', 'D19', 'D21', 'D22', 'D23', 'D25', 'D28', 'D29', 'D30', 'D32', 'D33', 'D38', 'D40', 'D41', 'D43', 'D45', 'D46', 'D47', 'D49'
)
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Price (D) and Volume/1h (E) updates.
$ws.Range('D2').Value = '26.638.31'
$ws.Range('E2').Value = '  +1.29%  '
$ws.Range('D3').Value = '1.632.84'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '213.10'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').Value = '0.500'
$ws.Range('E6').Value = '  +3.42%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +2.10%  '
$ws.Range('D9').Value = '0.0622'
$ws.Range('E9').Value = '  +1.24%  '
$ws.Range('D10').Value = '19.16'
$ws.Range('E10').Value = '  +2.10%  '
$ws.Range('D11').Value = '0.0843'
$ws.Range('E11').Value = '  +3.37%  '
$ws.Range('D12').Value = '1.862.33'
$ws.Range('E12').Value = '  +0.93%  '
$ws.Range('D13').Value = '1.624.76'
$ws.Range('E13').Value = '  +0.34%  '
$ws.Range('D14').Value = '4.09'
$ws.Range('E14').Value = '  +2.07%  '
$ws.Range('D15').Value = '0.525'
$ws.Range('E15').Value = '  +1.43%  '
$ws.Range('D16').Value = '26.636.51'
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('D17').Value = '63.30'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D18').Value = '0.0₃0741'
$ws.Range('E18').Value = '  +1.90%  '
$ws.Range('D19').Value = '219.47'
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').Value = '4.29'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').Value = '9.44'
$ws.Range('E22').Value = '  +1.31%  '
$ws.Range('D23').Value = '6.20'
$ws.Range('E23').Value = '  +2.62%  '
$ws.Range('E24').Value = '  +2.36%  '
$ws.Range('D25').Value = '148.75'
$ws.Range('E25').Value = '  +2.85%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('D28').Value = '6.92'
$ws.Range('E28').Value = '  +5.58%  '
$ws.Range('D29').Value = '15.51'
$ws.Range('E29').Value = '  +2.40%  '
$ws.Range('D30').Value = '0.0507'
$ws.Range('E30').Value = '  -2.86%  '
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('D32').Value = '3.29'
$ws.Range('E32').Value = '  +3.67%  '
$ws.Range('D33').Value = '2.98'
$ws.Range('E33').Value = '  +2.01%  '
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').Value = '1.210.43'
$ws.Range('E36').Value = '  +2.60%  '
$ws.Range('E37').Value = '  +5.84%  '
$ws.Range('D38').Value = '0.810'
$ws.Range('E38').Value = '  +0.49%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = '0.503'
$ws.Range('E40').Value = '  +1.58%  '
$ws.Range('D41').Value = '2.28'
$ws.Range('E41').Value = '  -1.37%  '
$ws.Range('E42').Value = '  +1.38%  '
$ws.Range('D43').Value = '0.792'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('D44').Value = '1.770.61'
$ws.Range('E44').Value = '  +0.81%  '
$ws.Range('D45').Value = '93.00'
$ws.Range('E45').Value = '  +0.44%  '
$ws.Range('D46').Value = '1.55'
$ws.Range('E46').Value = '  +1.16%  '
$ws.Range('D47').Value = '54.68'
$ws.Range('E47').Value = '  +1.75%  '
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('D49').Value = '7.68'
$ws.Range('E49').Value = '  +5.40%  '
$ws.Range('E50').Value = '  +0.37%  '
$ws.Range('E51').Value = '  +0.22%  '

# Restore default styling on the cells we temporarily forced to Text.
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).Style = "Normal"
}
